$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# Fix typo: MICHELE MATTIDORF -> MICHELE MATTIDORFF (row 17, column A)
$ws.Range("A17").Value = "MICHELE MATTIDORFF"

# Update the id next to EDUARDA SANTOS (row 2, column B)
$ws.Range("B2").Value = "5eaab222c733400015fa33d8"

# Move selection / view to match the author's final cursor position
$ws.Range("E23").Select()
$aw = $excel.ActiveWindow
$aw.ScrollRow = 13
$aw.ScrollColumn = 1
